# "antes de testar retirar apenas os 10 pixeis centrais para a sampple"
#
# Fill in the previously-empty sample row (row 9 - group 2 / "Castanho")
# with the measured HSV sample range, so the MINIFS/MAXIFS summary table
# (Q:W) picks it up on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 170
$ws.Range("D9").Value = 199
$ws.Range("E9").Value = 0.8
$ws.Range("F9").Value = 0.925
$ws.Range("G9").Value = 0.34
$ws.Range("H9").Value = 0.41

# Clear the stray border that was left on V3:V7 (Vmin column of the
# summary table) so it matches the rest of the column's formatting.
$ws.Range("V3:V7").Borders.LineStyle = -4142

# Leave the selection where the user ended up after the edit.
$ws.Range("C12").Select()
